# Generate Report for Handoff
# Adds two newly-handed-off files (5a306eb0-... and ac4d99b3-...) as rows 4 & 5
# on the "Overview" sheet and on each of the per-language status sheets
# ("zh-cn" and "de-de"), mirroring the existing rows for the other files.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# File identifiers involved in this handoff
# ---------------------------------------------------------------------------
$guid1 = "5a306eb0-0d0e-423e-9945-9caac84a27ec"
$hash1 = "019f5bba7f595988f05e926220d63d10563d25cc"
$guid2 = "ac4d99b3-93c4-4fb8-af98-0021f0923214"
$hash2 = "09ada6c7597e5394b450cea83e3000137010a4bb"

$status      = "Ready for handoff"
$extension   = ".md"
$reason      = "Include"
$noHandback  = "0001-01-01 00:00:00"

$overviewDate = "2016-03-24 22:41:05"
$zhDate       = "2016-03-24 22:40:59"
$deDate       = $overviewDate

# ===========================================================================
# Sheet "Overview": File Name | zh-cn | de-de | Latest Handoff Date
# ===========================================================================
$wsOverview = $wb.Worksheets.Item("Overview")

# --- row 4 : guid1 -----------------------------------------------------
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), `
    "https://github.com/OpenLocalizationTest/oltest/blob/072f605315fc9c1e4ef64890c6a8a1a5652771e9/e2e/$guid1.md", `
    "", "", "$guid1.md")
$wsOverview.Range("B4").Value = $status
$wsOverview.Range("C4").Value = $status
$wsOverview.Range("D4").Value = $overviewDate

# --- row 5 : guid2 -----------------------------------------------------
$wsOverview.Hyperlinks.Add($wsOverview.Range("A5"), `
    "https://github.com/OpenLocalizationTest/oltest/blob/c893e65a2d310ebcec1e2c38875afc24c86f1395/e2e/$guid2.md", `
    "", "", "$guid2.md")
$wsOverview.Range("B5").Value = $status
$wsOverview.Range("C5").Value = $status
$wsOverview.Range("D5").Value = $overviewDate

# ===========================================================================
# Sheet "zh-cn"
# ===========================================================================
$wsZh = $wb.Worksheets.Item("zh-cn")

# --- row 4 : guid1 -----------------------------------------------------
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), `
    "https://github.com/OpenLocalizationTest/oltest/blob/4e40eca9709ef02be553c724cde58d88411da8c1/e2e/$guid1.md", `
    "", "", "$guid1.md")
$wsZh.Range("B4").Value = $extension
$wsZh.Range("C4").Value = $status
$wsZh.Hyperlinks.Add($wsZh.Range("D4"), `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/21ddd4306ca5f3e07907428bace05e949f07bac7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$guid1.$hash1.zh-cn.xlf", `
    "", "", "$guid1.$hash1.zh-cn.xlf")
$wsZh.Range("E4").Value = $zhDate
$wsZh.Range("H4").Value = $noHandback
$wsZh.Range("J4").Value = $reason

# --- row 5 : guid2 -----------------------------------------------------
$wsZh.Hyperlinks.Add($wsZh.Range("A5"), `
    "https://github.com/OpenLocalizationTest/oltest/blob/545b8bb2b99db85c6ecb72716aa6833d622ec64a/e2e/$guid2.md", `
    "", "", "$guid2.md")
$wsZh.Range("B5").Value = $extension
$wsZh.Range("C5").Value = $status
$wsZh.Hyperlinks.Add($wsZh.Range("D5"), `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/27901c8c31cf1a436c574ec05b53cc6fd3044efe/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$guid2.$hash2.zh-cn.xlf", `
    "", "", "$guid2.$hash2.zh-cn.xlf")
$wsZh.Range("E5").Value = $zhDate
$wsZh.Range("H5").Value = $noHandback
$wsZh.Range("J5").Value = $reason

# ===========================================================================
# Sheet "de-de"
# ===========================================================================
$wsDe = $wb.Worksheets.Item("de-de")

# --- row 4 : guid1 -----------------------------------------------------
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), `
    "https://github.com/OpenLocalizationTest/oltest/blob/f96c2676554cad891d56a5f1a3513ac5a41b3b73/e2e/$guid1.md", `
    "", "", "$guid1.md")
$wsDe.Range("B4").Value = $extension
$wsDe.Range("C4").Value = $status
$wsDe.Hyperlinks.Add($wsDe.Range("D4"), `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c39f6fff7574c821ce113a287d24438fed2e65fa/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$guid1.$hash1.de-de.xlf", `
    "", "", "$guid1.$hash1.de-de.xlf")
$wsDe.Range("E4").Value = $deDate
$wsDe.Range("H4").Value = $noHandback
$wsDe.Range("J4").Value = $reason

# --- row 5 : guid2 -----------------------------------------------------
$wsDe.Hyperlinks.Add($wsDe.Range("A5"), `
    "https://github.com/OpenLocalizationTest/oltest/blob/f73ae19fa180a40ceebc263895c5d0cd84f42d63/e2e/$guid2.md", `
    "", "", "$guid2.md")
$wsDe.Range("B5").Value = $extension
$wsDe.Range("C5").Value = $status
$wsDe.Hyperlinks.Add($wsDe.Range("D5"), `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3e654e17d056e56827c3fe194ed6c97e95cda04e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$guid2.$hash2.de-de.xlf", `
    "", "", "$guid2.$hash2.de-de.xlf")
$wsDe.Range("E5").Value = $deDate
$wsDe.Range("H5").Value = $noHandback
$wsDe.Range("J5").Value = $reason
